# Update "想去人数" (number of people interested) counts for several
# events across the 展览, 演出 and 全部类型 sheets, as produced by the
# latest data refresh (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 79
$ws1.Range("F7").Value  = 96
$ws1.Range("F8").Value  = 262
$ws1.Range("F9").Value  = 19
$ws1.Range("F12").Value = 113
$ws1.Range("F13").Value = 2411
$ws1.Range("F14").Value = 30
$ws1.Range("F16").Value = 9
$ws1.Range("F18").Value = 528
$ws1.Range("F19").Value = 562
$ws1.Range("F20").Value = 165
$ws1.Range("F21").Value = 88
$ws1.Range("F24").Value = 1959
$ws1.Range("F25").Value = 4095
$ws1.Range("F28").Value = 1196
$ws1.Range("F30").Value = 2102
$ws1.Range("F34").Value = 121
$ws1.Range("F36").Value = 423
$ws1.Range("F38").Value = 707
$ws1.Range("F39").Value = 5
$ws1.Range("F40").Value = 440
$ws1.Range("F41").Value = 421

# ---- Sheet: 演出 ----
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 38

# ---- Sheet: 全部类型 ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 79
$ws4.Range("F7").Value  = 96
$ws4.Range("F8").Value  = 262
$ws4.Range("F9").Value  = 19
$ws4.Range("F12").Value = 113
$ws4.Range("F13").Value = 2411
$ws4.Range("F14").Value = 30
$ws4.Range("F16").Value = 38
$ws4.Range("F17").Value = 9
$ws4.Range("F19").Value = 528
$ws4.Range("F20").Value = 562
$ws4.Range("F21").Value = 165
$ws4.Range("F22").Value = 88
$ws4.Range("F25").Value = 1959
$ws4.Range("F26").Value = 4095
$ws4.Range("F29").Value = 1196
$ws4.Range("F31").Value = 2102
$ws4.Range("F35").Value = 121
$ws4.Range("F37").Value = 423
$ws4.Range("F39").Value = 707
$ws4.Range("F40").Value = 5
$ws4.Range("F41").Value = 440
$ws4.Range("F42").Value = 421
